$wb = $excel.ActiveWorkbook

$wsNeed = $wb.Worksheets.Item("Need")
$wsForm = $wb.Worksheets.Item("1. Formations")
$wsDev  = $wb.Worksheets.Item("2. Dev")

# --- Sheet "2. Dev": mark existing rows with an "X"/"x" in column A ---
$wsDev.Range("A8").Value = "X"
$wsDev.Range("A11").Value = "x"
$wsDev.Range("A12").Value = "x"

# --- New content, written in the exact order the strings were authored so ---
# --- that the shared-string table indices line up with the target file.  ---

# 1) sheet "2. Dev" row 13
$wsDev.Range("A13").Value = "X"
$wsDev.Range("B13").Value = 3.9
$wsDev.Range("C13").Value = "Size of chart needs to be bigger 3S-714 example plan view"

# 2) sheet "1. Formations" row 4
$wsForm.Range("B4").Value = 1.3
$wsForm.Range("C4").Value = "Update how formations are added (dynamically like pad view)"

# 3) sheet "2. Dev" row 14 (note the 2-decimal number format on B14)
$wsDev.Range("A14").Value = "X"
$wsDev.Range("B14").Value = 3.1
$wsDev.Range("B14").NumberFormat = "0.00"
$wsDev.Range("C14").Value = "uncheck box for formations if none entered"

# 4) sheet "2. Dev" row 15
$wsDev.Range("B15").Value = 3.11
$wsDev.Range("C15").Value = "Casing show depths on charts?"

# 5) sheet "2. Dev" row 16
$wsDev.Range("A16").Value = "X"
$wsDev.Range("B16").Value = 3.12
$wsDev.Range("C16").Value = "Formations names on side of the chart opposite well (if well N/S last point is - put them on the left, if + on right)"

# 6) sheet "Need" row 24
$wsNeed.Range("B24").Value = 1.5
$wsNeed.Range("C24").Value = "Generate DB file when none exists"

# 7) sheet "2. Dev" row 17
$wsDev.Range("A17").Value = "X"
$wsDev.Range("B17").Value = 3.13
$wsDev.Range("C17").Value = "close all charts when closing pad or chart view window"

# 8) sheet "2. Dev" row 18
$wsDev.Range("A18").Value = "X"
$wsDev.Range("B18").Value = 3.14
$wsDev.Range("C18").Value = "if show is not selected do not highlight"

# 9) sheet "2. Dev" row 19
$wsDev.Range("B19").Value = 3.15
$wsDev.Range("C19").Value = "annotate wells"

# 10) sheet "2. Dev" row 20
$wsDev.Range("A20").Value = "X"
$wsDev.Range("B20").Value = 3.16
$wsDev.Range("C20").Value = "highlight current well"

# 11) sheet "Need" row 25
$wsNeed.Range("A25").Value = "X"
$wsNeed.Range("B25").Value = 1.6
$wsNeed.Range("C25").Value = "Close all sub windows when main app window closes (Directional editor closes chart view and all sub windows close)"

# 12) sheet "2. Dev" row 21
$wsDev.Range("B21").Value = 3.17
$wsDev.Range("C21").Value = "Add offset well button"

# --- Selections on the sheets that end up inactive ---
[void]$wsNeed.Range("A26").Select()
[void]$wsForm.Range("G11").Select()

# --- Make "2. Dev" the active sheet/tab, with its final selection ---
[void]$wsDev.Activate()
[void]$wsDev.Range("B22").Select()
